$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ C=0.2189033714348056; D=0.181563596258755; E=0.1538151639921121; F=1.396127015626156; G=0.7888758353158636; H=0.8876228856497974; I=0.9677609843999697; J=0.17658885258545; K=1.3720748851091; L=0.2055227137656512; M=0.4322521183767378; O=3.368680463128655 }
  3 = @{ C=0.2165431892581324; D=0.178981191911106; E=0.1541985383094993; F=1.409464753198712; G=0.7993974093865361; H=0.8971542511202841; I=0.9775430492751234; J=0.1781929182568351; K=1.219001405342397; L=0.206983429835395; M=0.4021128343716143; O=3.410656766872208 }
  4 = @{ C=0.2151670069206375; D=0.1774449886228737; E=0.1544862713373902; F=1.418445943268999; G=0.8064239995019662; H=0.9034224952255201; I=0.9841093030528008; J=0.1792420378963797; K=1.124699482024141; L=0.2079471897897811; M=0.3836044263131484; O=3.438492064420174 }
  5 = @{ C=0.2146246458322167; D=0.1768314853837722; E=0.1546167163778467; F=1.422304952041223; G=0.8094296522743036; H=0.9060814756696658; I=0.9869259158111099; J=0.1796857249079533; K=1.086194639792353; L=0.2083567690147827; M=0.3760620171435818; O=3.45035335112216 }
  6 = @{ C=0.2145357039270266; D=0.1767303719648794; E=0.1546391743013604; F=1.422957762851738; G=0.8099373266523173; H=0.9065293169279443; I=0.9874021164895055; J=0.179760375380436; K=1.079796434883093; L=0.2084257970593857; M=0.374809618968527; O=3.452354198032623 }
  7 = @{ C=0.2151596176440904; D=0.1774366639244533; E=0.1544879771126872; F=1.418497181052963; G=0.8064639589678322; H=0.9034579315276261; I=0.9841467187556496; J=0.1792479561553471; K=1.124180496208055; L=0.2079526453065057; M=0.3835027062645722; O=3.438649932134894 }
  8 = @{ C=0.2180744885423991; D=0.1806629818994452; E=0.1539365063408908; F=1.400561539909404; G=0.7923860544062578; H=0.8908229873873736; I=0.9710175998499047; J=0.1771286149909788; K=1.319362082923874; L=0.206012502821677; M=0.4218610243489351; O=3.382725845081922 }
  9 = @{ C=0.2243659617108307; D=0.187378292037323; E=0.153269110166125; F=1.371672116916486; G=0.7692789276576448; H=0.8693449079835034; I=0.9497156579295307; J=0.1734814492440915; K=1.699508284553303; L=0.2027374707224148; M=0.4970365062623188; O=3.289428094081828 }
  10 = @{ C=0.2293348874204355; D=0.1925448336140505; E=0.1530296039520884; F=1.354277299707952; G=0.7550537264850661; H=0.8555740686003261; I=0.9367755231315513; J=0.1711111160177676; K=1.977092011178286; L=0.2006527261532192; M=0.5522149919310806; O=3.230875517651882 }
  11 = @{ C=0.2316697891734236; D=0.1949449644574628; E=0.1529747854928338; F=1.347195760099559; G=0.7491816205754134; H=0.849745178574338; I=0.9314775094429706; J=0.1700997310764141; K=2.102976933282605; L=0.1997738139764778; M=0.5773006172552329; O=3.206411413065666 }
  12 = @{ C=0.2325645859757941; D=0.1958609182532314; E=0.1529617829178314; F=1.344633738530732; G=0.7470442990723569; H=0.8476005358310843; I=0.9295559416262549; J=0.1697263504065045; K=2.150587729481913; L=0.1994509564603604; M=0.5867971493707529; O=3.197460076752279 }
  13 = @{ C=0.2323714046445389; D=0.1956633379441115; E=0.1529642387577503; F=1.345180195645824; G=0.7475007676408865; H=0.8480596374358953; I=0.9299660189600374; J=0.1698063374173433; K=2.140336570018803; L=0.1995200466026823; M=0.5847520382029074; O=3.199373994025592 }
  14 = @{ C=0.2317431923415114; D=0.1950201792297008; E=0.1529735604843125; F=1.346982583767335; G=0.7490040508867324; H=0.8495674822858632; I=0.9313177236249928; J=0.1700688203053815; K=2.106895104148975; L=0.1997470526414844; M=0.5780819634497476; O=3.205668712242215 }
  15 = @{ C=0.2313597744248597; D=0.1946271447883277; E=0.152980279526389; F=1.348102175712057; G=0.7499361017145958; H=0.8504992373168108; I=0.9321567103157378; J=0.1702308497128957; K=2.086403453352204; L=0.1998873978834084; M=0.5739959599319917; O=3.2095651418811 }
  16 = @{ C=0.2291837873872282; D=0.1923889746156107; E=0.1530342733521799; F=1.354756829520014; G=0.7554495515009947; H=0.8559637645121896; I=0.9371336072997067; J=0.1711785575691049; K=1.968857035183362; L=0.2007115609733106; M=0.5505752234590062; O=3.232518039396751 }
  17 = @{ C=0.2278679068574547; D=0.1910286311524345; E=0.1530812416118899; F=1.359052228876834; G=0.7589854241192739; H=0.8594276271691896; I=0.9403375396610087; J=0.171777070987563; K=1.896644183168462; L=0.2012349321166695; M=0.5362029710034903; O=3.247155418250301 }
  18 = @{ C=0.2271180662950627; D=0.1902508946749606; E=0.1531133535688767; F=1.361601086537306; G=0.7610755354211278; H=0.8614609437323466; I=0.9422357482279011; J=0.1721276168208661; K=1.855072801443214; L=0.2015424993265924; M=0.5279350253798611; O=3.255778844071642 }
  19 = @{ C=0.2268653915906498; D=0.1899883759869994; E=0.1531251024172988; F=1.362477525449016; G=0.7617928874202917; H=0.8621564305428109; I=0.9428879611237022; J=0.1722473873966273; K=1.840991279160846; L=0.2016477597875195; M=0.5251354192225861; O=3.258733680488717 }
  20 = @{ C=0.22800725882108; D=0.1911729565664047; E=0.1530757144405577; F=1.358586876607674; G=0.7586031889528044; H=0.8590546502948442; I=0.9399907424461986; J=0.1717127065958923; K=1.904335160540313; L=0.2011785418488437; M=0.537733072017005; O=3.245576087549807 }
  21 = @{ C=0.2319274260460276; D=0.1952088992105274; E=0.1529706121983416; F=1.346449932342324; G=0.7485601560904129; H=0.8491228918868217; I=0.9309183967765051; J=0.1699914620167338; K=2.116719301697913; L=0.1996801051035; M=0.5800412085354765; O=3.203811311503586 }
  22 = @{ C=0.234551340765293; D=0.1978878251857168; E=0.1529471164726139; F=1.339214898800464; G=0.7424996611943442; H=0.842996978400663; I=0.9254826660756379; J=0.1689225323360635; K=2.255178758382101; L=0.1987588791546138; M=0.6076751155649731; O=3.178338387949154 }
  23 = @{ C=0.2331452792125361; D=0.1964542918840095; E=0.1529555306489705; F=1.343012568431107; G=0.7456881585115909; H=0.8462330895549286; I=0.9283386408100256; J=0.1694879192391525; K=2.181313007413337; L=0.1992452459397889; M=0.5929281375728124; O=3.19176685378477 }
  24 = @{ C=0.2279442369955547; D=0.1911076935128762; E=0.1530781973578819; F=1.358797015123265; G=0.7587758190275125; H=0.8592231426366439; I=0.9401473543331491; J=0.1717417856380052; K=1.900858241536753; L=0.2012040150982344; M=0.5370413293301368; O=3.246289454453589 }
  25 = @{ C=0.2226027852000101; D=0.1855204493391369; E=0.1534054964318798; F=1.378814863556244; G=0.7750474181728535; H=0.8748023290960205; I=0.9550024981062322; J=0.1744137285232554; K=1.596960189264848; L=0.2035669060201393; M=0.4767070216770932; O=3.312913696378146 }
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $ws.Range("$col$row").Value = $data[$row][$col]
  }
}
